$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 531 - this shifts the existing rows 531..614 down to 532..615
$ws.Rows(531).Insert()

# Populate the new row 531 with the new data point (columns that are constant
# across this block are copied from the row below; the varying columns get
# the new values described by the diff)
$ws.Range("A531").Value = 4
$ws.Range("B531").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C531").Value = 'Los Lagos'
$ws.Range("D531").Value = 45218
$ws.Range("E531").Value = 10
$ws.Range("F531").Value = 100112008
$ws.Range("G531").Value = 'Coliflor'
$ws.Range("H531").Value = 'Sin especificar'
$ws.Range("I531").Value = 'Primera'
$ws.Range("J531").Value = 750
$ws.Range("K531").Value = 1500
$ws.Range("L531").Value = 1500
$ws.Range("M531").Value = 1500
$ws.Range("N531").Value = '$/unidad'
$ws.Range("O531").Value = 'Región Metropolitana'
$ws.Range("P531").Value = 1500
$ws.Range("Q531").Value = 1
$ws.Range("R531").Value = 'Hortaliza'
